# Weekly update: insert a new "Ají" price record as row 38 on the
# Terminal Hortofrutícola Agro Chillán sheet, pushing the existing rows
# 38-59 down to 39-60 (dimension grows from A1:R59 to A1:R60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 38; all rows below shift
# down by one (old row 38 becomes row 39, ..., old row 59 becomes row 60).
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly record.
$ws.Range("A38").Value = 7
$ws.Range("B38").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C38").Value = "Ñuble"
$ws.Range("D38").Value = 44582
$ws.Range("E38").Value = 16
$ws.Range("F38").Value = 100112021
$ws.Range("G38").Value = "Ají"
$ws.Range("H38").Value = "Americana (o)"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 80
$ws.Range("K38").Value = 16000
$ws.Range("L38").Value = 17000
$ws.Range("M38").Value = 16500
$ws.Range("N38").Value = "$/caja 15 kilos"
$ws.Range("O38").Value = "Región del Maule"
$ws.Range("P38").Value = 1100
$ws.Range("Q38").Value = 15
$ws.Range("R38").Value = "Hortaliza"
